# Update the workbook to reflect one additional day of carjacking data
# (rolling from "through December 03" to "through December 04" for 2022),
# which also revises the historical "December" columns for prior years
# as the underlying dataset is reprocessed to the same cutoff day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-12-04"

# Update the column header label (B1) for the current "through" month
$ws.Range("B1").Value = "December 2022 (through December 04)"

# Garfield Park - December 2020
$ws.Range("Z2").Value = 3

# Humboldt Park - December 2019
$ws.Range("AL3").Value = 1

# Roseland
$ws.Range("B5").Value = 1
$ws.Range("BJ5").Value = 1

# Woodlawn
$ws.Range("AX6").Value = 1
$ws.Range("N6").Value = 2

# Grand Crossing
$ws.Range("AL10").Value = 1
$ws.Range("N10").Value = 3

# New City
$ws.Range("N13").Value = 1

# Austin
$ws.Range("B15").Value = 1
$ws.Range("BV15").Value = 1
$ws.Range("N15").Value = 3

# North Lawndale
$ws.Range("Z19").Value = 1

# Avondale
$ws.Range("N23").Value = 2

# Ukrainian Village
$ws.Range("N24").Value = 1

# Portage Park
$ws.Range("AX26").Value = 2

# Auburn Gresham
$ws.Range("B28").Value = 1
$ws.Range("BJ28").Value = 2
$ws.Range("CH28").Value = 1

# Irving Park
$ws.Range("BJ31").Value = 2
$ws.Range("Z31").Value = 1

# Wicker Park
$ws.Range("Z34").Value = 1

# Gage Park
$ws.Range("BV35").Value = 1

# Calumet Heights
$ws.Range("N40").Value = 1

# Bridgeport
$ws.Range("BJ45").Value = 1
$ws.Range("Z45").Value = 1

# Chinatown
$ws.Range("N57").Value = 1

# Douglas
$ws.Range("BJ59").Value = 2
$ws.Range("N59").Value = 3

# Near South Side
$ws.Range("N83").Value = 1

# Sauganash,Forest Glen
$ws.Range("Z93").Value = 1
